$wb = $excel.ActiveWorkbook

# Move the "总计" (totals) sheet so it becomes the first sheet in the
# workbook, ahead of "2020-Q4".
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Move($wb.Worksheets.Item(1))

# "2020-Q4" stays the active/selected sheet (only the tab order changed).
$wb.Worksheets.Item("2020-Q4").Activate()
